# Applies the "Updated with ordnance usage from D12.1" change to the
# "D9 ->" worksheet (sheet1 in the package, first sheet in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New column header "Split to CVN" in L1:L2 (merged) -------------------
# Build the L1 format by copying the neighboring plain header cell (F1),
# then applying the centered / rotated 90 degree alignment used by the
# other single-column headers.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("L1").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("L1").HorizontalAlignment = -4108
$ws.Range("L1").Orientation = 90
$ws.Range("L1").Value2 = "Split to CVN"

# Build the L2 format by copying a cell that already carries the bottom
# thin border (B2), then tone the font back down to the small 10pt weight
# used throughout the header band, and restore default (bottom) vertical
# alignment so only the rotation/centering remain explicit.
$ws.Range("B2").Copy() | Out-Null
$ws.Range("L2").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false
$ws.Range("L2").Font.Bold = $false
$ws.Range("L2").Font.Size = 10
$ws.Range("L2").VerticalAlignment = -4107
$ws.Range("L2").Orientation = 90

$ws.Range("L1:L2").Merge()

# --- Shift the D-day sub headers in row 3 one column to the right ---------
# (L3 becomes the blank header above the new "Split to CVN" column, and the
# trailing "D16.2" header is retired.)
$ws.Range("L3").ClearContents()
$ws.Range("M3").Value2 = "D12.1"
$ws.Range("N3").Value2 = "D12.2"
$ws.Range("O3").Value2 = "D13.1"
$ws.Range("P3").Value2 = "D13.2"
$ws.Range("Q3").Value2 = "D14.1"
$ws.Range("R3").Value2 = "D14.2"
$ws.Range("S3").Value2 = "D15.1"
$ws.Range("T3").Value2 = "D15.2"
$ws.Range("U3").Value2 = "D16.1"

# --- Advance the "current as of" marker from D12.1 to D12.2 ---------------
$ws.Range("V1").Value2 = "D12.2"

# --- Ordnance usage recorded against D12.1 (now column M) -----------------
$ws.Range("M13").Value2 = 5
$ws.Range("M18").Value2 = 25

# --- Selection bookmark ----------------------------------------------------
$ws.Range("Q7").Select() | Out-Null
